# Laboratorio 7 - Entrega final
# Update the measured PROBING / CHAINING data with the final benchmark
# results, and leave the selection where it was left in the authoring
# session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- PROBING table (rows 3:5) -------------------------------------------
$ws.Range("B3").Value = 1327329.4129999999
$ws.Range("C3").Value = 45145.758999999998

$ws.Range("B4").Value = 1327329.4369999999
$ws.Range("C4").Value = 48539.574000000001

$ws.Range("B5").Value = 1327329.4369999999
$ws.Range("C5").Value = 44607.85

# --- CHAINING table (rows 10:12) ----------------------------------------
$ws.Range("B10").Value = 1327340.2949999999
$ws.Range("C10").Value = 47278.074999999997

$ws.Range("B11").Value = 1327341.443
$ws.Range("C11").Value = 45541.419000000002

$ws.Range("B12").Value = 1327341.443
$ws.Range("C12").Value = 47953.781999999999

# --- Final cursor position left by the author ----------------------------
$ws.Range("C12").Select() | Out-Null
